$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values.
# All Price cells are stored as text in the source data (e.g. "23.977.61" or
# "0.3840"). Plain decimal strings would otherwise be auto-converted to numbers
# by Excel (dropping significant trailing zeros / changing the stored type), so
# force a Text number format on those cells before writing the new value.

$ws.Range("D2").Value = "23.977.61"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.655.06"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.84"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3901"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3840"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.37"
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.355"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08454"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.93"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.129"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.901"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001316"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").Value = "1.654.14"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.64"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06996"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.78"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.942"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.71"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "23.974.01"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.482"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.004"
$ws.Range("E26").Value = "  +6.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.14"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.44"
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.461"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.36"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.850"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "1.835.95"
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.046"
$ws.Range("E34").Value = "  +6.77%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02972"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.770"
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.92"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2689"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09157"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7562"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.48"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.430"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.37"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6955"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.459"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.092"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.84"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.219"
$ws.Range("E51").Value = "  +1.21%  "
